$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 15: continuation of the effort log table, copying the date
# formatting (style) of the cell above so no new number format is created.
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A15").Value2 = 41443
$ws.Range("B15").Value2 = 1.5
$ws.Range("C15").Value2 = 2.5
$ws.Range("D15").Value2 = $ws.Range("D14").Value2

$ws.Range("A15").Select()
